$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add Sheet2 right after Sheet1 -----------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# --- Populate data. Assignment order below controls the shared-strings ----
# --- table append order so it reproduces the target workbook exactly. -----

# Row 1 (headers)
$ws2.Range("A1").Value = "TestName"
$ws2.Range("B1").Value = "TestCode"
$ws2.Range("C1").Value = "WorkFiled"
$ws2.Range("D1").Value = "Discription"
$ws2.Range("E1").Value = "Show on Customer Portal"
$ws2.Range("F1").Value = "Dynamic"

# Column A (names) rows 2-9
$ws2.Range("A2").Value = "Raghu"
$ws2.Range("A3").Value = "Timmy"
$ws2.Range("A4").Value = "Taj Mohan"
$ws2.Range("A5").Value = "Noobmaster"
$ws2.Range("A6").Value = "Pika"
$ws2.Range("A7").Value = "DD"
$ws2.Range("A8").Value = "Nec"
$ws2.Range("A9").Value = "Der"

# Column C (work field) - note non-sequential discovery order
$ws2.Range("C5").Value = "ACCOUNT_PAYABLE"
$ws2.Range("C4").Value = "QUOTE"
$ws2.Range("C6").Value = "CONSOLIDATION"
$ws2.Range("C7").Value = "CUSTOMER_BOOKING"
$ws2.Range("C8").Value = "SHIPMENT"
$ws2.Range("C9").Value = "CUSTOMER_INTEGRATION"

# Column D (discription) rows 2-9
$ws2.Range("D2").Value = "New msg"
$ws2.Range("D3").Value = "milestone"
$ws2.Range("D4").Value = "Hello milestone"
$ws2.Range("D5").Value = "wow milestone"
$ws2.Range("D6").Value = "rock milestone"
$ws2.Range("D7").Value = "create mielstone"
$ws2.Range("D8").Value = "fine with"
$ws2.Range("D9").Value = "ok done"

# Remaining column C cells that reuse already-existing shared strings
$ws2.Range("C2").Value = "BOOKING"
$ws2.Range("C3").Value = "ORDER"

# Column B (numbers)
$ws2.Range("B2").Value = 15678
$ws2.Range("B3").Value = 6701
$ws2.Range("B4").Value = 9078
$ws2.Range("B5").Value = 1672
$ws2.Range("B6").Value = 4599
$ws2.Range("B7").Value = 6022
$ws2.Range("B8").Value = 9033
$ws2.Range("B9").Value = 2311

# Column E (Show on Customer Portal - boolean)
$ws2.Range("E2").Value = $true
$ws2.Range("E3").Value = $false
$ws2.Range("E4").Value = $true
$ws2.Range("E5").Value = $false
$ws2.Range("E6").Value = $true
$ws2.Range("E7").Value = $false
$ws2.Range("E8").Value = $true
$ws2.Range("E9").Value = $false

# Column F (Dynamic - boolean)
$ws2.Range("F2").Value = $false
$ws2.Range("F3").Value = $true
$ws2.Range("F4").Value = $true
$ws2.Range("F5").Value = $false
$ws2.Range("F6").Value = $false
$ws2.Range("F7").Value = $true
$ws2.Range("F8").Value = $true
$ws2.Range("F9").Value = $false

# --- Formatting: QUOTE cell and the 4 "false" Customer-Portal cells get a --
# --- 9pt Courier New font -----------------------------------------------
$ws2.Range("C4").Font.Size = 9
$ws2.Range("C4").Font.Name = "Courier New"

$ws2.Range("C4").Copy()
$ws2.Range("E3").PasteSpecial(-4122)
$ws2.Range("E5").PasteSpecial(-4122)
$ws2.Range("E7").PasteSpecial(-4122)
$ws2.Range("E9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column widths (approximate Excel's AutoFit result) -------------------
$ws2.Columns.Item(1).ColumnWidth = 11.084
$ws2.Columns.Item(2).ColumnWidth = 8.417
$ws2.Columns.Item(3).ColumnWidth = 23.417
$ws2.Columns.Item(4).ColumnWidth = 15.251
$ws2.Columns.Item(5).ColumnWidth = 22.751
$ws2.Columns.Item(6).ColumnWidth = 7.584

# --- Selections -------------------------------------------------------------
$ws2.Activate()
$ws2.Range("F11").Select()

$ws1.Activate()
$ws1.Range("D19").Select()
